$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1459
$ws.Range("F3").Value = 1429
$ws.Range("F6").Value = 689
$ws.Range("F8").Value = 626
$ws.Range("F9").Value = 472
$ws.Range("F11").Value = 1374
$ws.Range("F12").Value = 32517
$ws.Range("F13").Value = 6983
$ws.Range("F15").Value = 356
$ws.Range("F16").Value = 570
$ws.Range("F17").Value = 435
$ws.Range("F19").Value = 98
$ws.Range("F20").Value = 14
$ws.Range("F21").Value = 44
$ws.Range("F22").Value = 442
$ws.Range("F23").Value = 101
$ws.Range("F24").Value = 788
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = 316
$ws.Range("F27").Value = 387
$ws.Range("F28").Value = 436
$ws.Range("F30").Value = 191
$ws.Range("F31").Value = 47
$ws.Range("F32").Value = 735
$ws.Range("F33").Value = 288
$ws.Range("F34").Value = 132
$ws.Range("F35").Value = 728
$ws.Range("F36").Value = 108
$ws.Range("F38").Value = 788
$ws.Range("F40").Value = 53
$ws.Range("F41").Value = 22

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 1172
$ws.Range("F3").Value = 7
$ws.Range("F5").Value = 155
$ws.Range("F19").Value = 4291

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1440
$ws.Range("F3").Value = 351

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1440
$ws.Range("F3").Value = 351
$ws.Range("F4").Value = 1172
$ws.Range("F5").Value = 1459
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 1429
$ws.Range("F9").Value = 689
$ws.Range("F11").Value = 626
$ws.Range("F13").Value = 1374
$ws.Range("F14").Value = 155
$ws.Range("F21").Value = 6984
$ws.Range("F23").Value = 356
$ws.Range("F25").Value = 570
$ws.Range("F26").Value = 435
$ws.Range("F28").Value = 98
$ws.Range("F31").Value = 442
$ws.Range("F32").Value = 101
$ws.Range("F33").Value = 788
$ws.Range("F34").Value = 316
$ws.Range("F35").Value = 387
$ws.Range("F36").Value = 436
$ws.Range("F38").Value = 191
$ws.Range("F39").Value = 47
$ws.Range("F40").Value = 735
$ws.Range("F42").Value = 288
$ws.Range("F43").Value = 132
$ws.Range("F44").Value = 108
$ws.Range("F45").Value = 788
$ws.Range("F47").Value = 53
$ws.Range("F49").Value = 22
